$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D15").PasteSpecial(-4122)

$ws.Range("E15").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E15").PasteSpecial(-4122)

$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("F15").PasteSpecial(-4122)

$ws.Range("H15").Value = -100

$ws.Range("M15").Value = -50

$ws.Range("N15").Value = -62.5

$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C16").PasteSpecial(-4122)

$ws.Range("E16").Value = -100

$ws.Range("G16").Value = 3

$ws.Range("H16").Value = 66.666666666666

$ws.Range("J16").Value = 37

$ws.Range("K16").Value = -13.513513513513

$ws.Range("M16").Value = -48.387096774193

$ws.Range("N16").Value = -80.722891566265

$ws.Range("C17").Value = 5

$ws.Range("D17").Value = 4

$ws.Range("E17").Value = 25

$ws.Range("F17").Value = 13

$ws.Range("G17").Value = 17

$ws.Range("H17").Value = -23.529411764705

$ws.Range("I17").Value = 112

$ws.Range("J17").Value = 65

$ws.Range("K17").Value = 72.307692307692

$ws.Range("L17").Value = 53.424657534246

$ws.Range("M17").Value = 34.939759036144

$ws.Range("N17").Value = -37.430167597765

$ws.Range("C18").Value = 3
$ws.Range("I14").Copy()
$ws.Range("C18").PasteSpecial(-4122)

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D18").PasteSpecial(-4122)

$ws.Range("E18").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E18").PasteSpecial(-4122)

$ws.Range("F18").Value = 11

$ws.Range("G18").Value = 4

$ws.Range("H18").Value = 175

$ws.Range("I18").Value = 61

$ws.Range("K18").Value = 56.410256410256

$ws.Range("L18").Value = 29.787234042553

$ws.Range("M18").Value = -47.413793103448

$ws.Range("N18").Value = -92.422360248447

$ws.Range("C19").Value = 7

$ws.Range("D19").Value = 7

$ws.Range("E19").Value = 0

$ws.Range("F19").Value = 43

$ws.Range("G19").Value = 33

$ws.Range("H19").Value = 30.30303030303

$ws.Range("I19").Value = 282

$ws.Range("J19").Value = 198

$ws.Range("K19").Value = 42.424242424242

$ws.Range("L19").Value = 39.603960396039

$ws.Range("M19").Value = 16.528925619834

$ws.Range("N19").Value = -43.6

$ws.Range("D20").Value = 6

$ws.Range("E20").Value = -33.333333333333

$ws.Range("F20").Value = 18

$ws.Range("G20").Value = 13

$ws.Range("H20").Value = 38.461538461538

$ws.Range("I20").Value = 68

$ws.Range("J20").Value = 78

$ws.Range("K20").Value = -12.820512820512

$ws.Range("L20").Value = 142.857142857143

$ws.Range("M20").Value = 0

$ws.Range("N20").Value = -95.930580490724

$ws.Range("C21").Value = 19

$ws.Range("E21").Value = 5.555555555555

$ws.Range("F21").Value = 90

$ws.Range("G21").Value = 71

$ws.Range("H21").Value = 26.760563380281

$ws.Range("I21").Value = 562

$ws.Range("J21").Value = 423

$ws.Range("K21").Value = 32.860520094562

$ws.Range("L21").Value = 47.120418848167

$ws.Range("M21").Value = -3.602058319039

$ws.Range("N21").Value = -83.17365269461

$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C23").PasteSpecial(-4122)

$ws.Range("D23").Value = 1
$ws.Range("I14").Copy()
$ws.Range("D23").PasteSpecial(-4122)

$ws.Range("E23").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E23").PasteSpecial(-4122)

$ws.Range("F23").Value = 1

$ws.Range("G23").Value = 1
$ws.Range("I14").Copy()
$ws.Range("G23").PasteSpecial(-4122)

$ws.Range("H23").Value = 0
$ws.Range("K14").Copy()
$ws.Range("H23").PasteSpecial(-4122)

$ws.Range("J23").Value = 13

$ws.Range("K23").Value = 69.230769230769

$ws.Range("M23").Value = 83.333333333333

$ws.Range("C24").Value = 28

$ws.Range("D24").Value = 19

$ws.Range("E24").Value = 47.368421052631

$ws.Range("F24").Value = 89

$ws.Range("G24").Value = 90

$ws.Range("H24").Value = -1.111111111111

$ws.Range("I24").Value = 679

$ws.Range("J24").Value = 451

$ws.Range("K24").Value = 50.554323725055

$ws.Range("L24").Value = 142.5

$ws.Range("M24").Value = -34.774255523535

$ws.Range("C25").Value = 7

$ws.Range("D25").Value = 3

$ws.Range("E25").Value = 133.333333333333

$ws.Range("F25").Value = 35

$ws.Range("G25").Value = 25

$ws.Range("H25").Value = 40

$ws.Range("I25").Value = 193

$ws.Range("J25").Value = 208

$ws.Range("K25").Value = -7.211538461538

$ws.Range("L25").Value = 34.965034965035

$ws.Range("M25").Value = -43.401759530791

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D26").PasteSpecial(-4122)

$ws.Range("E26").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E26").PasteSpecial(-4122)

$ws.Range("F26").NumberFormat = "@"
$ws.Range("F26").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("F26").PasteSpecial(-4122)

$ws.Range("H26").Value = -100

$ws.Range("C27").Value = 1
$ws.Range("I14").Copy()
$ws.Range("C27").PasteSpecial(-4122)

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D27").PasteSpecial(-4122)

$ws.Range("E27").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E27").PasteSpecial(-4122)

$ws.Range("F27").Value = 3

$ws.Range("H27").Value = 50

$ws.Range("I27").Value = 23

$ws.Range("K27").Value = 9.523809523809

$ws.Range("L27").Value = 27.777777777777

$ws.Range("C28").Value = 1
$ws.Range("I14").Copy()
$ws.Range("C28").PasteSpecial(-4122)

$ws.Range("F28").Value = 1
$ws.Range("I14").Copy()
$ws.Range("F28").PasteSpecial(-4122)

$ws.Range("I28").Value = 2

$ws.Range("M28").Value = 100

$ws.Range("N28").Value = 0

$ws.Range("C29").Value = 1
$ws.Range("I14").Copy()
$ws.Range("C29").PasteSpecial(-4122)

$ws.Range("F29").Value = 1
$ws.Range("I14").Copy()
$ws.Range("F29").PasteSpecial(-4122)

$ws.Range("I29").Value = 2

$ws.Range("M29").Value = 100

$ws.Range("N29").Value = 0

$ws.Range("C30").Value = 1
$ws.Range("I14").Copy()
$ws.Range("C30").PasteSpecial(-4122)

$ws.Range("F30").Value = 1
$ws.Range("I14").Copy()
$ws.Range("F30").PasteSpecial(-4122)

$ws.Range("H30").Value = 0

$ws.Range("I30").Value = 5

$ws.Range("K30").Value = -16.666666666666

$ws.Range("L30").Value = 66.666666666666

$ws.Range("A8").Value = "Volume 30   Number  31"
$ws.Range("C9").Value = "Report Covering the Week  7/31/2023  Through  8/6/2023"

$excel.CutCopyMode = 0